$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "short-url" column (B2:B62) from "xxBC7k" to "6RIyAV" ---
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "6RIyAV"
}

# --- Update "refugees" value for row 61 (N61) from 16 to 19 ---
# Use a scratch cell + text formula + PasteSpecial(values-only) so the
# numeric-looking text "19" lands as TEXT (matching the existing column
# content) instead of being auto-coerced into a Number, while keeping the
# destination cell's existing formatting/style untouched.
$ws.Range("ZZ1").Formula = "=""19"""
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("N61").PasteSpecial(-4163) | Out-Null
$ws.Range("ZZ1").Clear() | Out-Null

# --- Update "stateless" value for row 62 (S62) from 26811 to 20590 ---
$ws.Range("ZZ1").Formula = "=""20590"""
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("S62").PasteSpecial(-4163) | Out-Null
$ws.Range("ZZ1").Clear() | Out-Null

$excel.CutCopyMode = 0
